$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 22657
$ws.Range("J3").Value = 22657
$ws.Range("L3").Value = 22657
$ws.Range("N3").Value = -22885
$ws.Range("H33").Value = 414.5
$ws.Range("I33").Value = 434.18182
$ws.Range("K33").Value = 434.18182
$ws.Range("M33").Value = -205.18182
$ws.Range("H39").Value = 160.5
$ws.Range("I39").Value = 43.4
$ws.Range("K39").Value = 130.2
$ws.Range("M39").Value = 165.8
$ws.Range("H40").Value = 2069.077
$ws.Range("I40").Value = 1300.3334
$ws.Range("J40").Value = 2728
$ws.Range("K40").Value = 1300.3334
$ws.Range("L40").Value = 2728
$ws.Range("M40").Value = -1125.3334
$ws.Range("N40").Value = -3078
$ws.Range("H69").Value = 10002.5
$ws.Range("I69").Value = 8005
$ws.Range("J69").Value = 12000
$ws.Range("K69").Value = 24015
$ws.Range("L69").Value = 36000
$ws.Range("M69").Value = -23141
$ws.Range("N69").Value = -37748
$ws.Range("H72").Value = 10002.5
$ws.Range("I72").Value = 8005
$ws.Range("J72").Value = 12000
$ws.Range("K72").Value = 72045
$ws.Range("L72").Value = 108000
$ws.Range("M72").Value = -67677
$ws.Range("N72").Value = -116736
$ws.Range("H92").Value = 801.3125
$ws.Range("I92").Value = 780.3
$ws.Range("J92").Value = 836.3333
$ws.Range("K92").Value = 780.3
$ws.Range("L92").Value = 836.3333
$ws.Range("M92").Value = 467.7
$ws.Range("N92").Value = -3332.3333
$ws.Range("H99").Value = 173.25
$ws.Range("I99").Value = 170.66667
$ws.Range("J99").Value = 181
$ws.Range("K99").Value = 512.00001
$ws.Range("L99").Value = 543
$ws.Range("M99").Value = 985.99999
$ws.Range("N99").Value = -3539
$ws.Range("H101").Value = 395
$ws.Range("I101").Value = 395
$ws.Range("K101").Value = 1185
$ws.Range("M101").Value = 437
$ws.Range("H102").Value = 22657
$ws.Range("J102").Value = 22657
$ws.Range("L102").Value = 22657
$ws.Range("N102").Value = -29147
$ws.Range("H106").Value = 9208.450000000001
$ws.Range("I106").Value = 6676.222
$ws.Range("K106").Value = 6676.222
$ws.Range("M106").Value = -6045.222
$ws.Range("H118").Value = 340
$ws.Range("I118").Value = 340
$ws.Range("K118").Value = 1020
$ws.Range("M118").Value = 637
$ws.Range("H125").Value = 583.3333
$ws.Range("J125").Value = 500
$ws.Range("L125").Value = 4500
$ws.Range("N125").Value = -9420
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("H131").Value = 689.36365
$ws.Range("I131").Value = 735.625
$ws.Range("J131").Value = 566
$ws.Range("K131").Value = 2206.875
$ws.Range("L131").Value = 1698
$ws.Range("M131").Value = 2833.125
$ws.Range("N131").Value = -11778
$ws.Range("H132").Value = 4100.3477
$ws.Range("I132").Value = 4100.3477
$ws.Range("K132").Value = 12301.0431
$ws.Range("M132").Value = -9771.043100000001
$ws.Range("H137").Value = 4634.143
$ws.Range("J137").Value = 6185.75
$ws.Range("L137").Value = 18557.25
$ws.Range("N137").Value = -23657.25
$ws.Range("H138").Value = 2822.2222
$ws.Range("I138").Value = 1080.4
$ws.Range("K138").Value = 3241.2
$ws.Range("M138").Value = 1898.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 379.75
$ws.Range("I5").Value = 139
$ws.Range("K5").Value = 139
$ws.Range("M5").Value = -27
$ws.Range("H13").Value = 2511499.8
$ws.Range("J13").Value = 15333
$ws.Range("L13").Value = 15333
$ws.Range("N13").Value = -15621
$ws.Range("H30").Value = 3752.3333
$ws.Range("I30").Value = 838
$ws.Range("J30").Value = 6666.6665
$ws.Range("K30").Value = 838
$ws.Range("L30").Value = 6666.6665
$ws.Range("M30").Value = -688
$ws.Range("N30").Value = -6966.6665
$ws.Range("H32").Value = 35286
$ws.Range("I32").Value = 35286
$ws.Range("K32").Value = 35286
$ws.Range("M32").Value = -34999
$ws.Range("H45").Value = 1481.7142
$ws.Range("I45").Value = 1449.5385
$ws.Range("K45").Value = 1449.5385
$ws.Range("M45").Value = -1072.5385
$ws.Range("H61").Value = 1529.2858
$ws.Range("I61").Value = 1529.2858
$ws.Range("K61").Value = 1529.2858
$ws.Range("M61").Value = -1317.2858
$ws.Range("H92").Value = 95277
$ws.Range("J92").Value = 95277
$ws.Range("L92").Value = 95277
$ws.Range("N92").Value = -100269
$ws.Range("H97").Value = 843659.0600000001
$ws.Range("I97").Value = 1123652.1
$ws.Range("J97").Value = 3680
$ws.Range("K97").Value = 1123652.1
$ws.Range("L97").Value = 3680
$ws.Range("M97").Value = -1123156.1
$ws.Range("N97").Value = -4672
$ws.Range("H110").Value = 11360
$ws.Range("I110").Value = 11360
$ws.Range("K110").Value = 11360
$ws.Range("M110").Value = -9315
$ws.Range("H122").Value = 2344.85
$ws.Range("I122").Value = 2344.85
$ws.Range("K122").Value = 7034.549999999999
$ws.Range("M122").Value = -4584.549999999999
$ws.Range("H136").Value = 1529.2858
$ws.Range("I136").Value = 1529.2858
$ws.Range("K136").Value = 4587.857400000001
$ws.Range("M136").Value = -2037.857400000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 379.75
$ws.Range("I4").Value = 139
$ws.Range("K4").Value = 139
$ws.Range("M4").Value = -24
$ws.Range("H20").Value = 3480.4
$ws.Range("I20").Value = 2469.3333
$ws.Range("J20").Value = 4997
$ws.Range("K20").Value = 2469.3333
$ws.Range("L20").Value = 4997
$ws.Range("M20").Value = -2222.3333
$ws.Range("N20").Value = -5491
$ws.Range("H22").Value = 493.5
$ws.Range("I22").Value = 491.66666
$ws.Range("K22").Value = 491.66666
$ws.Range("M22").Value = -318.66666
$ws.Range("H86").Value = 3251.9285
$ws.Range("I86").Value = 3425.2307
$ws.Range("J86").Value = 999
$ws.Range("K86").Value = 3425.2307
$ws.Range("L86").Value = 999
$ws.Range("M86").Value = -2302.2307
$ws.Range("N86").Value = -3245
$ws.Range("H89").Value = 3251.9285
$ws.Range("I89").Value = 3425.2307
$ws.Range("J89").Value = 999
$ws.Range("K89").Value = 17126.1535
$ws.Range("L89").Value = 4995
$ws.Range("M89").Value = -11510.1535
$ws.Range("N89").Value = -16227
$ws.Range("H105").Value = 4031.889
$ws.Range("I105").Value = 2321
$ws.Range("K105").Value = 2321
$ws.Range("M105").Value = -574
$ws.Range("H130").Value = 90000
$ws.Range("J130").Value = 90000
$ws.Range("L130").Value = 90000
$ws.Range("N130").Value = -100040
$ws.Range("H131").Value = 52500
$ws.Range("I131").Value = 52500
$ws.Range("K131").Value = 52500
$ws.Range("M131").Value = -47460

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1500
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H16").Value = 1737
$ws.Range("I16").Value = 1737
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1737
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1450
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 9499.866
$ws.Range("I22").Value = 1856.4286
$ws.Range("J22").Value = 16187.875
$ws.Range("K22").Value = 1856.4286
$ws.Range("L22").Value = 16187.875
$ws.Range("M22").Value = -1506.4286
$ws.Range("N22").Value = -16887.875
$ws.Range("H26").Value = 4850
$ws.Range("J26").Value = 5000
$ws.Range("L26").Value = 5000
$ws.Range("N26").Value = -5574
$ws.Range("H31").Value = 2733.2856
$ws.Range("I31").Value = 3047.2
$ws.Range("J31").Value = 1948.5
$ws.Range("K31").Value = 3047.2
$ws.Range("L31").Value = 1948.5
$ws.Range("M31").Value = -2752.2
$ws.Range("N31").Value = -2538.5
$ws.Range("H34").Value = 2733.2856
$ws.Range("I34").Value = 3047.2
$ws.Range("J34").Value = 1948.5
$ws.Range("K34").Value = 3047.2
$ws.Range("L34").Value = 1948.5
$ws.Range("M34").Value = -2845.2
$ws.Range("N34").Value = -2352.5
$ws.Range("H62").Value = 2876.6
$ws.Range("I62").Value = 2747.5
$ws.Range("K62").Value = 2747.5
$ws.Range("M62").Value = -2123.5
$ws.Range("H65").Value = 2876.6
$ws.Range("I65").Value = 2747.5
$ws.Range("K65").Value = 13737.5
$ws.Range("M65").Value = -10617.5
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502
$ws.Range("H105").Value = 5000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 5000
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -8494
$ws.Range("H113").Value = 1737
$ws.Range("I113").Value = 1737
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1737
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 433
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
$ws.Range("H132").Value = 2511.2222
$ws.Range("I132").Value = 2511.2222
$ws.Range("K132").Value = 7533.6666
$ws.Range("M132").Value = -5003.6666

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 243.5
$ws.Range("J2").Value = 351.36365
$ws.Range("L2").Value = 2108.1819
$ws.Range("N2").Value = -2334.1819
$ws.Range("H7").Value = 19181.125
$ws.Range("I7").Value = 19181.125
$ws.Range("K7").Value = 57543.375
$ws.Range("M7").Value = -57431.375
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H112").Value = 10411.3
$ws.Range("I112").Value = 2056.5
$ws.Range("J112").Value = 12500
$ws.Range("K112").Value = 6169.5
$ws.Range("L112").Value = 37500
$ws.Range("M112").Value = -5061.5
$ws.Range("N112").Value = -39716
$ws.Range("H117").Value = 436.75
$ws.Range("J117").Value = 449
$ws.Range("L117").Value = 1347
$ws.Range("N117").Value = -8231
$ws.Range("H126").Value = 1030
$ws.Range("I126").Value = 1030
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3090
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 1850
$ws.Range("N126").ClearContents()
$ws.Range("H129").Value = 1924.4
$ws.Range("I129").Value = 1263.2
$ws.Range("J129").Value = 2585.6
$ws.Range("K129").Value = 3789.6
$ws.Range("L129").Value = 7756.799999999999
$ws.Range("M129").Value = 1210.4
$ws.Range("N129").Value = -17756.8

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 153749.88
$ws.Range("I11").Value = 138333.33
$ws.Range("K11").Value = 138333.33
$ws.Range("M11").Value = -138194.33
$ws.Range("H21").Value = 49500
$ws.Range("I21").Value = 49000
$ws.Range("J21").Value = 50000
$ws.Range("K21").Value = 49000
$ws.Range("L21").Value = 50000
$ws.Range("M21").Value = -48827
$ws.Range("N21").Value = -50346
$ws.Range("H30").Value = 49500
$ws.Range("I30").Value = 49000
$ws.Range("J30").Value = 50000
$ws.Range("K30").Value = 49000
$ws.Range("L30").Value = 50000
$ws.Range("M30").Value = -48895
$ws.Range("N30").Value = -50210
$ws.Range("H34").Value = 25000
$ws.Range("J34").Value = 25000
$ws.Range("L34").Value = 25000
$ws.Range("N34").Value = -25536
$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25630
$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27184
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -54900
$ws.Range("H126").Value = 1145
$ws.Range("I126").Value = 1145
$ws.Range("K126").Value = 3435
$ws.Range("M126").Value = -965
$ws.Range("H132").Value = 1999.6666
$ws.Range("I132").Value = 1999.6666
$ws.Range("K132").Value = 5998.9998
$ws.Range("M132").Value = -3468.9998
$ws.Range("H141").Value = 97666.664
$ws.Range("J141").Value = 97666.664
$ws.Range("L141").Value = 97666.664
$ws.Range("N141").Value = -108026.664

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H16").Value = 565.1875
$ws.Range("I16").Value = 367.42856
$ws.Range("K16").Value = 367.42856
$ws.Range("M16").Value = -197.42856
$ws.Range("H22").Value = 2350
$ws.Range("I22").Value = 1800
$ws.Range("K22").Value = 1800
$ws.Range("M22").Value = -1505
$ws.Range("H23").Value = 4999.5
$ws.Range("I23").Value = 4999.5
$ws.Range("K23").Value = 4999.5
$ws.Range("M23").Value = -4769.5
$ws.Range("H26").Value = 349.5
$ws.Range("I26").Value = 399
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 399
$ws.Range("L26").Value = 300
$ws.Range("M26").Value = -104
$ws.Range("N26").Value = -890
$ws.Range("H27").Value = 2350
$ws.Range("I27").Value = 1800
$ws.Range("K27").Value = 1800
$ws.Range("M27").Value = -1693
$ws.Range("H46").Value = 3891
$ws.Range("I46").Value = 1400
$ws.Range("K46").Value = 1400
$ws.Range("M46").Value = -1212
$ws.Range("H61").Value = 4497.5
$ws.Range("I61").Value = 4497.5
$ws.Range("K61").Value = 4497.5
$ws.Range("M61").Value = -4295.5
$ws.Range("H111").Value = 19998
$ws.Range("J111").Value = 19998
$ws.Range("L111").Value = 19998
$ws.Range("N111").Value = -28178
$ws.Range("H113").Value = 4497.5
$ws.Range("I113").Value = 4497.5
$ws.Range("K113").Value = 4497.5
$ws.Range("M113").Value = -2327.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 12521124
$ws.Range("I2").Value = 16688165
$ws.Range("J2").Value = 19999
$ws.Range("K2").Value = 16688165
$ws.Range("L2").Value = 19999
$ws.Range("M2").Value = -16688053
$ws.Range("N2").Value = -20223
$ws.Range("H25").Value = 25000
$ws.Range("J25").Value = 25000
$ws.Range("L25").Value = 25000
$ws.Range("N25").Value = -25586
$ws.Range("H26").Value = 19500
$ws.Range("I26").Value = 5000
$ws.Range("J26").Value = 34000
$ws.Range("K26").Value = 5000
$ws.Range("L26").Value = 34000
$ws.Range("M26").Value = -4707
$ws.Range("N26").Value = -34586
$ws.Range("H29").Value = 48005
$ws.Range("J29").Value = 90000
$ws.Range("L29").Value = 90000
$ws.Range("N29").Value = -90580
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H100").Value = 1248
$ws.Range("I100").Value = 330.66666
$ws.Range("K100").Value = 661.33332
$ws.Range("M100").Value = -120.33332
